$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet previously had a stray leftover cell B72 (just whitespace) after
# row 71 ("game_lost-connection-modal_btn-txt"). We now add three new
# localization rows (72-74) for the "Get Free Chips" modal, reusing the
# same visual styling already used by the table (fill = style 1, fill+wrap
# = style 2) so no new cell formats are introduced.
# ---------------------------------------------------------------------------

# --- Row 72: global_get-free-chips-modal_header -----------------------------
# Style pattern for this row: A=plain fill, B=fill+wrap, C=plain fill
$ws.Range("A71").Copy()
$ws.Range("A72").PasteSpecial(-4122)
$ws.Range("C72").PasteSpecial(-4122)

$ws.Range("B71").Copy()
$ws.Range("B72").PasteSpecial(-4122)

$ws.Cells.Item(72, 1).Value = "global_get-free-chips-modal_header"
$ws.Cells.Item(72, 2).Value = "Refuel  "
$ws.Cells.Item(72, 3).Value = "Auftanken"

# --- Row 73: global_get-free-chips-modal_content -----------------------------
# Style pattern for this row matches row 71: A=plain fill, B=fill+wrap, C=fill+wrap
$ws.Range("A71:C71").Copy()
$ws.Range("A73:C73").PasteSpecial(-4122)

$ws.Cells.Item(73, 1).Value = "global_get-free-chips-modal_content"
$ws.Cells.Item(73, 2).Value = "Oh noes, it seems like you're running out of chips! But don't worry, here's a fresh batch of chips for you so you can continue playing!"
$ws.Cells.Item(73, 3).Value = "Oh nein, es scheint, als würden Ihnen die Chips ausgehen! Aber keine Sorge, hier ist eine neue Charge Chips für Sie, damit Sie weiterspielen können!"

# --- Row 74: global_get-free-chips-modal_btn-txt -----------------------------
# Style pattern for this row: A=plain fill, B=plain fill, C=plain fill (no wrap)
$ws.Range("A70:C70").Copy()
$ws.Range("A74:C74").PasteSpecial(-4122)

$ws.Cells.Item(74, 1).Value = "global_get-free-chips-modal_btn-txt"
$ws.Cells.Item(74, 2).Value = "Get Your Free Chips"
$ws.Cells.Item(74, 3).Value = "Gratis Chips Holen"

$excel.CutCopyMode = 0
